# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 3 = "R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 477
$wsOff.Range("C3").Value = 334
$wsOff.Range("D3").Value = 108
$wsOff.Range("E3").Value = 56

# Update DEF sheet (row 3 = "R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 465
$wsDef.Range("C3").Value = 325
$wsDef.Range("D3").Value = 111
$wsDef.Range("E3").Value = 56
